$wb = $excel.ActiveWorkbook
$wsClassroom = $wb.Worksheets.Item("Classroom_Allocation")

$wsClassroom.Range("G6").Value = "large classroom"
$wsClassroom.Range("H6").NumberFormat = "@"
$wsClassroom.Range("H6").Value = "120"
$wsClassroom.Range("H6").Style = "Normal"
$wsClassroom.Range("I6").Value = ""
$wsClassroom.Range("M6").Value = "C001"

$wsClassroom.Range("I7").Value = "Projector"
$wsClassroom.Range("M7").Value = "C101"

$wsClassroom.Range("I8").Value = "Projector"
$wsClassroom.Range("M8").Value = "C202"

$wsClassroom.Range("M9").Value = "C204"

$wsClassroom.Range("G17").Value = "classroom"
$wsClassroom.Range("H17").NumberFormat = "@"
$wsClassroom.Range("H17").Value = "96"
$wsClassroom.Range("H17").Style = "Normal"
$wsClassroom.Range("I17").Value = "Projector"
$wsClassroom.Range("M17").Value = "C102"

$wsClassroom.Range("G18").Value = "classroom"
$wsClassroom.Range("H18").NumberFormat = "@"
$wsClassroom.Range("H18").Value = "96"
$wsClassroom.Range("H18").Style = "Normal"
$wsClassroom.Range("I18").Value = "Projector"
$wsClassroom.Range("M18").Value = "C104"

$wsClassroom.Range("G22").Value = "large classroom"
$wsClassroom.Range("H22").NumberFormat = "@"
$wsClassroom.Range("H22").Value = "120"
$wsClassroom.Range("H22").Style = "Normal"
$wsClassroom.Range("I22").Value = ""
$wsClassroom.Range("M22").Value = "C001"

$wsClassroom.Range("I23").Value = "Projector"
$wsClassroom.Range("M23").Value = "C101"

$wsClassroom.Range("I24").Value = "Projector"
$wsClassroom.Range("M24").Value = "C202"

$wsClassroom.Range("M25").Value = "C204"

$wsClassroom.Range("G39").Value = "classroom"
$wsClassroom.Range("H39").NumberFormat = "@"
$wsClassroom.Range("H39").Value = "96"
$wsClassroom.Range("H39").Style = "Normal"
$wsClassroom.Range("M39").Value = "C101"

$wsClassroom.Range("M40").Value = "C202"

$wsClassroom.Range("M41").Value = "C204"

$wsClassroom.Range("G42").Value = "classroom"
$wsClassroom.Range("H42").NumberFormat = "@"
$wsClassroom.Range("H42").Value = "96"
$wsClassroom.Range("H42").Style = "Normal"
$wsClassroom.Range("I42").Value = "Projector"
$wsClassroom.Range("M42").Value = "C102"

$wsClassroom.Range("G43").Value = "classroom"
$wsClassroom.Range("H43").NumberFormat = "@"
$wsClassroom.Range("H43").Value = "96"
$wsClassroom.Range("H43").Style = "Normal"
$wsClassroom.Range("I43").Value = "Projector"
$wsClassroom.Range("M43").Value = "C104"

$wsClassroom.Range("G46").Value = "classroom"
$wsClassroom.Range("H46").NumberFormat = "@"
$wsClassroom.Range("H46").Value = "96"
$wsClassroom.Range("H46").Style = "Normal"
$wsClassroom.Range("I46").Value = "Projector"
$wsClassroom.Range("M46").Value = "C102"

$wsClassroom.Range("G47").Value = "classroom"
$wsClassroom.Range("H47").NumberFormat = "@"
$wsClassroom.Range("H47").Value = "96"
$wsClassroom.Range("H47").Style = "Normal"
$wsClassroom.Range("I47").Value = "Projector"
$wsClassroom.Range("M47").Value = "C104"

$wsClassroom.Range("G48").Value = "classroom"
$wsClassroom.Range("H48").NumberFormat = "@"
$wsClassroom.Range("H48").Value = "96"
$wsClassroom.Range("H48").Style = "Normal"
$wsClassroom.Range("M48").Value = "C202"

$wsClassroom.Range("I49").Value = "TV"
$wsClassroom.Range("M49").Value = "C203"

$wsClassroom.Range("I53").Value = "Projector"
$wsClassroom.Range("M53").Value = "C002"

$wsClassroom.Range("G54").Value = "classroom"
$wsClassroom.Range("H54").NumberFormat = "@"
$wsClassroom.Range("H54").Value = "96"
$wsClassroom.Range("H54").Style = "Normal"
$wsClassroom.Range("M54").Value = "C102"

$wsClassroom.Range("M55").Value = "C104"

$wsClassroom.Range("I56").Value = "TV"
$wsClassroom.Range("M56").Value = "C203"

$wsClassroom.Range("I57").Value = "TV"
$wsClassroom.Range("M57").Value = "C205"

$wsClassroom.Range("G65").Value = "Auditorium"
$wsClassroom.Range("H65").NumberFormat = "@"
$wsClassroom.Range("H65").Value = "240"
$wsClassroom.Range("H65").Style = "Normal"
$wsClassroom.Range("I65").Value = "Audio/Video System"
$wsClassroom.Range("M65").Value = "C004"

$wsClassroom.Range("G66").Value = "large classroom"
$wsClassroom.Range("H66").NumberFormat = "@"
$wsClassroom.Range("H66").Value = "120"
$wsClassroom.Range("H66").Style = "Normal"
$wsClassroom.Range("I66").Value = ""
$wsClassroom.Range("M66").Value = "C001"

$wsClassroom.Range("G67").Value = "large classroom"
$wsClassroom.Range("H67").NumberFormat = "@"
$wsClassroom.Range("H67").Value = "120"
$wsClassroom.Range("H67").Style = "Normal"
$wsClassroom.Range("M67").Value = "C002"

$wsClassroom.Range("M68").Value = "C101"

$wsClassroom.Range("I69").Value = "Projector"
$wsClassroom.Range("M69").Value = "C002"

$wsClassroom.Range("G70").Value = "classroom"
$wsClassroom.Range("H70").NumberFormat = "@"
$wsClassroom.Range("H70").Value = "96"
$wsClassroom.Range("H70").Style = "Normal"
$wsClassroom.Range("M70").Value = "C102"

$wsClassroom.Range("M71").Value = "C104"

$wsClassroom.Range("I72").Value = "TV"
$wsClassroom.Range("M72").Value = "C203"

$wsClassroom.Range("I73").Value = "TV"
$wsClassroom.Range("M73").Value = "C205"

$wsClassroom.Range("G85").Value = "large classroom"
$wsClassroom.Range("H85").NumberFormat = "@"
$wsClassroom.Range("H85").Value = "120"
$wsClassroom.Range("H85").Style = "Normal"
$wsClassroom.Range("M85").Value = "C002"

$wsClassroom.Range("M86").Value = "C102"

$wsClassroom.Range("M87").Value = "C104"

$wsClassroom.Range("M89").Value = "C205"

$wsClassroom.Range("G90").Value = "Auditorium"
$wsClassroom.Range("H90").NumberFormat = "@"
$wsClassroom.Range("H90").Value = "240"
$wsClassroom.Range("H90").Style = "Normal"
$wsClassroom.Range("I90").Value = "Audio/Video System"
$wsClassroom.Range("M90").Value = "C004"

$wsClassroom.Range("G91").Value = "large classroom"
$wsClassroom.Range("H91").NumberFormat = "@"
$wsClassroom.Range("H91").Value = "120"
$wsClassroom.Range("H91").Style = "Normal"
$wsClassroom.Range("I91").Value = ""
$wsClassroom.Range("M91").Value = "C001"

$wsClassroom.Range("G92").Value = "large classroom"
$wsClassroom.Range("H92").NumberFormat = "@"
$wsClassroom.Range("H92").Value = "120"
$wsClassroom.Range("H92").Style = "Normal"
$wsClassroom.Range("M92").Value = "C002"

$wsClassroom.Range("M93").Value = "C101"

$wsClassroom.Range("G94").Value = "Auditorium"
$wsClassroom.Range("H94").NumberFormat = "@"
$wsClassroom.Range("H94").Value = "240"
$wsClassroom.Range("H94").Style = "Normal"
$wsClassroom.Range("I94").Value = "Audio/Video System"
$wsClassroom.Range("M94").Value = "C004"

$wsClassroom.Range("G95").Value = "large classroom"
$wsClassroom.Range("H95").NumberFormat = "@"
$wsClassroom.Range("H95").Value = "120"
$wsClassroom.Range("H95").Style = "Normal"
$wsClassroom.Range("I95").Value = ""
$wsClassroom.Range("M95").Value = "C001"

$wsClassroom.Range("G96").Value = "large classroom"
$wsClassroom.Range("H96").NumberFormat = "@"
$wsClassroom.Range("H96").Value = "120"
$wsClassroom.Range("H96").Style = "Normal"
$wsClassroom.Range("M96").Value = "C002"

$wsClassroom.Range("I97").Value = "Projector"
$wsClassroom.Range("M97").Value = "C101"

$wsBasket = $wb.Worksheets.Item("Basket_Course_Allocations")

$wsBasket.Range("C9").Value = "C002, C004"
$wsBasket.Range("C10").Value = "C001, C102"
$wsBasket.Range("C11").Value = "C101, C104"
$wsBasket.Range("C12").Value = "C202, C203"
$wsBasket.Range("C13").Value = "C204, C205"
$wsBasket.Range("C14").Value = "C004, C102"
$wsBasket.Range("C15").Value = "C001, C104"
$wsBasket.Range("C16").Value = "C002, C202"
$wsBasket.Range("C17").Value = "C101, C203"

Write-Output "done"